# Fruta / hortaliza, semanal
# Insert one new weekly record at row 367 (pushing the existing rows
# 367-398 down to 368-399) on the only worksheet in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 367:398 down one row, leaving a blank row 367 in place.
$ws.Rows.Item(367).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(367, 1).Value  = 4
$ws.Cells.Item(367, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(367, 3).Value  = "Los Lagos"
$ws.Cells.Item(367, 4).Value  = 45106
$ws.Cells.Item(367, 5).Value  = 10
$ws.Cells.Item(367, 6).Value  = "Fruta"
$ws.Cells.Item(367, 7).Value  = 100108
$ws.Cells.Item(367, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(367, 9).Value  = 100108002
$ws.Cells.Item(367, 10).Value = "Mango"
$ws.Cells.Item(367, 11).Value = "Sin especificar"
$ws.Cells.Item(367, 12).Value = "Primera"
$ws.Cells.Item(367, 13).Value = 200
$ws.Cells.Item(367, 14).Value = 8500
$ws.Cells.Item(367, 15).Value = 9000
$ws.Cells.Item(367, 16).Value = 8750
$ws.Cells.Item(367, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(367, 18).Value = "Brasil"
$ws.Cells.Item(367, 19).Value = 2188
$ws.Cells.Item(367, 20).Value = 4
